$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3200
$ws.Range("I69").Value = 2800
$ws.Range("J69").Value = 4000
$ws.Range("K69").Value = 8400
$ws.Range("L69").Value = 12000
$ws.Range("M69").Value = -7526
$ws.Range("N69").Value = -13748

$ws.Range("H70").Value = 22395.8
$ws.Range("J70").Value = 27619.75
$ws.Range("L70").Value = 82859.25
$ws.Range("N70").Value = -83399.25

$ws.Range("H72").Value = 3200
$ws.Range("I72").Value = 2800
$ws.Range("J72").Value = 4000
$ws.Range("K72").Value = 25200
$ws.Range("L72").Value = 36000
$ws.Range("M72").Value = -20832
$ws.Range("N72").Value = -44736

$ws.Range("H73").Value = 22395.8
$ws.Range("J73").Value = 27619.75
$ws.Range("L73").Value = 82859.25
$ws.Range("N73").Value = -84731.25

$ws.Range("H86").Value = 1099.6
$ws.Range("I86").Value = 999.3333
$ws.Range("J86").Value = 1250
$ws.Range("K86").Value = 999.3333
$ws.Range("L86").Value = 1250
$ws.Range("M86").Value = 123.6667
$ws.Range("N86").Value = -3496

$ws.Range("H89").Value = 1099.6
$ws.Range("I89").Value = 999.3333
$ws.Range("J89").Value = 1250
$ws.Range("K89").Value = 4996.6665
$ws.Range("L89").Value = 6250
$ws.Range("M89").Value = 619.3334999999997
$ws.Range("N89").Value = -17482

$ws.Range("H132").Value = 1269.5
$ws.Range("I132").Value = 1138.4445
$ws.Range("K132").Value = 3415.3335
$ws.Range("M132").Value = -885.3335000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3305.7869
$ws.Range("I32").Value = 1807.6078
$ws.Range("K32").Value = 1807.6078
$ws.Range("M32").Value = -1520.6078

$ws.Range("H61").Value = 3388.0588
$ws.Range("I61").Value = 2456.4211
$ws.Range("K61").Value = 2456.4211
$ws.Range("M61").Value = -2244.4211

$ws.Range("H74").Value = 1716.4783
$ws.Range("I74").Value = 1215.6666
$ws.Range("J74").Value = 3519.4
$ws.Range("K74").Value = 1215.6666
$ws.Range("L74").Value = 3519.4
$ws.Range("M74").Value = -341.6666
$ws.Range("N74").Value = -5267.4

$ws.Range("H77").Value = 1716.4783
$ws.Range("I77").Value = 1215.6666
$ws.Range("J77").Value = 3519.4
$ws.Range("K77").Value = 6078.333000000001
$ws.Range("L77").Value = 17597
$ws.Range("M77").Value = -1710.333000000001
$ws.Range("N77").Value = -26333

$ws.Range("H88").Value = 4469.7
$ws.Range("J88").Value = 4985.2856
$ws.Range("L88").Value = 4985.2856
$ws.Range("N88").Value = -5797.2856

$ws.Range("H91").Value = 4469.7
$ws.Range("J91").Value = 4985.2856
$ws.Range("L91").Value = 4985.2856
$ws.Range("N91").Value = -7793.2856

$ws.Range("H97").Value = 1381.1428
$ws.Range("I97").Value = 1164.25
$ws.Range("K97").Value = 1164.25
$ws.Range("M97").Value = -668.25

$ws.Range("H102").Value = 2500
$ws.Range("I102").Value = 2250
$ws.Range("K102").Value = 2250
$ws.Range("M102").Value = -628

$ws.Range("H110").Value = 1306.6428
$ws.Range("I110").Value = 1007.03845
$ws.Range("J110").Value = 5201.5
$ws.Range("K110").Value = 1007.03845
$ws.Range("L110").Value = 5201.5
$ws.Range("M110").Value = 1037.96155
$ws.Range("N110").Value = -9291.5

$ws.Range("H132").Value = 1542.3043
$ws.Range("I132").Value = 1093.2162
$ws.Range("K132").Value = 3279.6486
$ws.Range("M132").Value = -749.6486000000004

$ws.Range("H136").Value = 3388.0588
$ws.Range("I136").Value = 2456.4211
$ws.Range("K136").Value = 7369.263300000001
$ws.Range("M136").Value = -4819.263300000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 228
$ws.Range("J22").Value = 310
$ws.Range("L22").Value = 310
$ws.Range("N22").Value = -656

$ws.Range("H86").Value = 127713.25
$ws.Range("J86").Value = 251470.12
$ws.Range("L86").Value = 251470.12
$ws.Range("N86").Value = -253716.12

$ws.Range("H89").Value = 127713.25
$ws.Range("J89").Value = 251470.12
$ws.Range("L89").Value = 1257350.6
$ws.Range("N89").Value = -1268582.6

$ws.Range("H94").Value = 1251.8
$ws.Range("I94").Value = 1064.5
$ws.Range("K94").Value = 1064.5
$ws.Range("M94").Value = -613.5

$ws.Range("H134").Value = 13334.083
$ws.Range("I134").Value = 19013.625
$ws.Range("K134").Value = 57040.875
$ws.Range("M134").Value = -54505.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1404134.4
$ws.Range("I58").Value = 1812715.4
$ws.Range("J58").Value = 3285.4285
$ws.Range("K58").Value = 1812715.4
$ws.Range("L58").Value = 3285.4285
$ws.Range("M58").Value = -1812512.4
$ws.Range("N58").Value = -3691.4285

$ws.Range("H132").Value = 2022.1316
$ws.Range("I132").Value = 1247.2693
$ws.Range("K132").Value = 3741.8079
$ws.Range("M132").Value = -1211.8079

$ws.Range("H136").Value = 1404134.4
$ws.Range("I136").Value = 1812715.4
$ws.Range("J136").Value = 3285.4285
$ws.Range("K136").Value = 5438146.199999999
$ws.Range("L136").Value = 9856.2855
$ws.Range("M136").Value = -5435596.199999999
$ws.Range("N136").Value = -14956.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 226.16667
$ws.Range("I2").Value = 318.33334
$ws.Range("K2").Value = 1910.00004
$ws.Range("M2").Value = -1797.00004

$ws.Range("H99").Value = 2970
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 2970
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 8910
$ws.Range("N99").Value = -13402
$ws.Range("M99").ClearContents()

$ws.Range("H128").Value = 397499.25
$ws.Range("I128").Value = 397499.25
$ws.Range("K128").Value = 1192497.75
$ws.Range("M128").Value = -1187517.75

$ws.Range("H131").Value = 1815.37
$ws.Range("I131").Value = 534
$ws.Range("J131").Value = 1973.7416
$ws.Range("K131").Value = 1602
$ws.Range("L131").Value = 5921.2248
$ws.Range("M131").Value = 3438
$ws.Range("N131").Value = -16001.2248

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3316.2727
$ws.Range("I102").Value = 3348
$ws.Range("J102").Value = 2999
$ws.Range("K102").Value = 3348
$ws.Range("L102").Value = 2999
$ws.Range("M102").Value = -1726
$ws.Range("N102").Value = -6243

$ws.Range("H122").Value = 1904.7333
$ws.Range("I122").Value = 1537.3
$ws.Range("J122").Value = 2639.6
$ws.Range("K122").Value = 4611.9
$ws.Range("L122").Value = 7918.799999999999
$ws.Range("M122").Value = -2161.9
$ws.Range("N122").Value = -12818.8

$ws.Range("H127").Value = 34557
$ws.Range("J127").Value = 34557
$ws.Range("L127").Value = 34557
$ws.Range("N127").Value = -44477

$ws.Range("H132").Value = 2960655.5
$ws.Range("I132").Value = 3498148.8
$ws.Range("K132").Value = 10494446.4
$ws.Range("M132").Value = -10491916.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3286.4119
$ws.Range("I16").Value = 4224.625
$ws.Range("K16").Value = 4224.625
$ws.Range("M16").Value = -4054.625

$ws.Range("H64").Value = 268972
$ws.Range("J64").Value = 25296.334
$ws.Range("L64").Value = 25296.334
$ws.Range("N64").Value = -25746.334

$ws.Range("H67").Value = 268972
$ws.Range("J67").Value = 25296.334
$ws.Range("L67").Value = 25296.334
$ws.Range("N67").Value = -26856.334

$ws.Range("H82").Value = 897
$ws.Range("I82").Value = 897
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 897
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -536
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 897
$ws.Range("I85").Value = 897
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 897
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 351
$ws.Range("N85").ClearContents()

$ws.Range("H93").Value = 575.5
$ws.Range("I93").Value = 530.6
$ws.Range("K93").Value = 530.6
$ws.Range("M93").Value = 717.4

$ws.Range("H132").Value = 2039.7872
$ws.Range("I132").Value = 1829.6111
$ws.Range("J132").Value = 2170.2415
$ws.Range("K132").Value = 5488.8333
$ws.Range("L132").Value = 6510.7245
$ws.Range("M132").Value = -2958.8333
$ws.Range("N132").Value = -11570.7245

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 12629.667
$ws.Range("J18").Value = 15000
$ws.Range("L18").Value = 15000
$ws.Range("N18").Value = -15346

$ws.Range("H63").Value = 22999.5
$ws.Range("J63").Value = 22999.5
$ws.Range("L63").Value = 22999.5
$ws.Range("N63").Value = -24247.5

$ws.Range("H66").Value = 22999.5
$ws.Range("J66").Value = 22999.5
$ws.Range("L66").Value = 68998.5
$ws.Range("N66").Value = -75238.5

$ws.Range("H81").Value = 1833.5
$ws.Range("I81").Value = 1500.25
$ws.Range("K81").Value = 3000.5
$ws.Range("M81").Value = -1939.5

$ws.Range("H84").Value = 1833.5
$ws.Range("I84").Value = 1500.25
$ws.Range("K84").Value = 15002.5
$ws.Range("M84").Value = -9698.5

$ws.Range("H126").Value = 3192.1765
$ws.Range("I126").Value = 1783.6154
$ws.Range("K126").Value = 5350.8462
$ws.Range("M126").Value = -2880.8462

$ws.Range("H132").Value = 1486.8462
$ws.Range("I132").Value = 1143.381
$ws.Range("K132").Value = 3430.143
$ws.Range("M132").Value = -900.143
